$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 532.35
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 532.35
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1597.05
$ws.Range("N17").Value = -1933.05
$ws.Range("M17").Value = $null

$ws.Range("H40").Value = 2383.875
$ws.Range("J40").Value = 2967.75
$ws.Range("L40").Value = 2967.75
$ws.Range("N40").Value = -3317.75

$ws.Range("H62").Value = 5252.5625
$ws.Range("I62").Value = 5951.222
$ws.Range("J62").Value = 4354.2856
$ws.Range("K62").Value = 5951.222
$ws.Range("L62").Value = 4354.2856
$ws.Range("M62").Value = -5327.222
$ws.Range("N62").Value = -5602.2856

$ws.Range("H65").Value = 5252.5625
$ws.Range("I65").Value = 5951.222
$ws.Range("J65").Value = 4354.2856
$ws.Range("K65").Value = 29756.11
$ws.Range("L65").Value = 21771.428
$ws.Range("M65").Value = -26636.11
$ws.Range("N65").Value = -28011.428

$ws.Range("H76").Value = 4550.0713
$ws.Range("I76").Value = 3903
$ws.Range("J76").Value = 4599.846
$ws.Range("K76").Value = 3903
$ws.Range("L76").Value = 4599.846
$ws.Range("M76").Value = -3588
$ws.Range("N76").Value = -5229.846

$ws.Range("H79").Value = 4550.0713
$ws.Range("I79").Value = 3903
$ws.Range("J79").Value = 4599.846
$ws.Range("K79").Value = 3903
$ws.Range("L79").Value = 4599.846
$ws.Range("M79").Value = -2811
$ws.Range("N79").Value = -6783.846

$ws.Range("H101").Value = 830.4
$ws.Range("J101").Value = 1061.3334
$ws.Range("L101").Value = 3184.0002
$ws.Range("N101").Value = -6428.0002

$ws.Range("H135").Value = 2171.3225
$ws.Range("I135").Value = 1784
$ws.Range("J135").Value = 3499.2856
$ws.Range("K135").Value = 16056
$ws.Range("L135").Value = 31493.5704
$ws.Range("M135").Value = -13521
$ws.Range("N135").Value = -36563.5704

$ws.Range("H137").Value = 32689.55
$ws.Range("I137").Value = 23286.857
$ws.Range("K137").Value = 69860.571
$ws.Range("M137").Value = -67310.571

$ws.Range("H138").Value = 20676.873
$ws.Range("J138").Value = 38941.5
$ws.Range("L138").Value = 116824.5
$ws.Range("N138").Value = -127104.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20556.219
$ws.Range("I32").Value = 21435.441
$ws.Range("J32").Value = 5316.3335
$ws.Range("K32").Value = 21435.441
$ws.Range("L32").Value = 5316.3335
$ws.Range("M32").Value = -21148.441
$ws.Range("N32").Value = -5890.3335

$ws.Range("H102").Value = 4982.25
$ws.Range("I102").Value = 3409.5
$ws.Range("K102").Value = 3409.5
$ws.Range("M102").Value = -1787.5

$ws.Range("H110").Value = 36556.92
$ws.Range("I110").Value = 39059.477
$ws.Range("K110").Value = 39059.477
$ws.Range("M110").Value = -37014.477

$ws.Range("H132").Value = 2853.8462
$ws.Range("I132").Value = 2414.3809
$ws.Range("K132").Value = 7243.1427
$ws.Range("M132").Value = -4713.1427

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2247.5
$ws.Range("I86").Value = 1996
$ws.Range("J86").Value = 2624.75
$ws.Range("K86").Value = 1996
$ws.Range("L86").Value = 2624.75
$ws.Range("M86").Value = -873
$ws.Range("N86").Value = -4870.75

$ws.Range("H89").Value = 2247.5
$ws.Range("I89").Value = 1996
$ws.Range("J89").Value = 2624.75
$ws.Range("K89").Value = 9980
$ws.Range("L89").Value = 13123.75
$ws.Range("M89").Value = -4364
$ws.Range("N89").Value = -24355.75

$ws.Range("H94").Value = 4930.2173
$ws.Range("I94").Value = 7083.2144
$ws.Range("K94").Value = 7083.2144
$ws.Range("M94").Value = -6632.2144

$ws.Range("H105").Value = 1614.8334
$ws.Range("I105").Value = 1454.7307
$ws.Range("J105").Value = 2655.5
$ws.Range("K105").Value = 1454.7307
$ws.Range("L105").Value = 2655.5
$ws.Range("M105").Value = 292.2692999999999
$ws.Range("N105").Value = -6149.5

$ws.Range("H107").Value = 1893.7941
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null

$ws.Range("H134").Value = 2664.1292
$ws.Range("I134").Value = 2455.7144
$ws.Range("K134").Value = 7367.1432
$ws.Range("M134").Value = -4832.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5267822.5
$ws.Range("I31").Value = 10004033
$ws.Range("J31").Value = 5365.8887
$ws.Range("K31").Value = 10004033
$ws.Range("L31").Value = 5365.8887
$ws.Range("M31").Value = -10003738
$ws.Range("N31").Value = -5955.8887

$ws.Range("H34").Value = 5267822.5
$ws.Range("I34").Value = 10004033
$ws.Range("J34").Value = 5365.8887
$ws.Range("K34").Value = 10004033
$ws.Range("L34").Value = 5365.8887
$ws.Range("M34").Value = -10003831
$ws.Range("N34").Value = -5769.8887

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null

$ws.Range("H58").Value = 1513.2122
$ws.Range("I58").Value = 1349.6538
$ws.Range("J58").Value = 2120.7144
$ws.Range("K58").Value = 1349.6538
$ws.Range("L58").Value = 2120.7144
$ws.Range("M58").Value = -1146.6538
$ws.Range("N58").Value = -2526.7144

$ws.Range("H60").Value = 30999.5
$ws.Range("J60").Value = 30999.5
$ws.Range("L60").Value = 30999.5
$ws.Range("N60").Value = -32021.5

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null

$ws.Range("H135").Value = 120000
$ws.Range("J135").Value = 120000
$ws.Range("L135").Value = 120000
$ws.Range("M135").Value = -130140

$ws.Range("H136").Value = 1513.2122
$ws.Range("I136").Value = 1349.6538
$ws.Range("J136").Value = 2120.7144
$ws.Range("K136").Value = 4048.9614
$ws.Range("L136").Value = 6362.1432
$ws.Range("M136").Value = -1498.9614
$ws.Range("N136").Value = -11462.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4983.3335
$ws.Range("J80").Value = 4983.3335
$ws.Range("L80").Value = 14950.0005
$ws.Range("N80").Value = -16822.0005

$ws.Range("H83").Value = 4983.3335
$ws.Range("J83").Value = 4983.3335
$ws.Range("L83").Value = 44850.0015
$ws.Range("N83").Value = -54210.0015

$ws.Range("H86").Value = 724.3333
$ws.Range("I86").Value = 724.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2172.9999
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -986.9998999999998
$ws.Range("N86").Value = $null

$ws.Range("H89").Value = 724.3333
$ws.Range("I89").Value = 724.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6518.9997
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -590.9997000000003
$ws.Range("N89").Value = $null

$ws.Range("H107").Value = 1830.8125
$ws.Range("J107").Value = 930.6
$ws.Range("L107").Value = 2791.8
$ws.Range("N107").Value = -6631.8

$ws.Range("H129").Value = 2755.111
$ws.Range("J129").Value = 2983.6
$ws.Range("L129").Value = 8950.799999999999
$ws.Range("N129").Value = -18950.8

$ws.Range("H139").Value = 5632.2
$ws.Range("I139").Value = 5632.2
$ws.Range("K139").Value = 16896.6
$ws.Range("M139").Value = -11756.6

$ws.Range("H141").Value = 6283.4287
$ws.Range("J141").Value = 8333.333000000001
$ws.Range("L141").Value = 24999.999
$ws.Range("N141").Value = -35359.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18049.5
$ws.Range("J15").Value = 18666.111
$ws.Range("L15").Value = 18666.111
$ws.Range("N15").Value = -19242.111

$ws.Range("H70").Value = 6104.143
$ws.Range("I70").Value = 5247
$ws.Range("J70").Value = 6447
$ws.Range("K70").Value = 5247
$ws.Range("L70").Value = 6447
$ws.Range("M70").Value = -4977
$ws.Range("N70").Value = -6987

$ws.Range("H73").Value = 6104.143
$ws.Range("I73").Value = 5247
$ws.Range("J73").Value = 6447
$ws.Range("K73").Value = 5247
$ws.Range("L73").Value = 6447
$ws.Range("M73").Value = -4311
$ws.Range("N73").Value = -8319

$ws.Range("H81").Value = 18049.5
$ws.Range("J81").Value = 18666.111
$ws.Range("L81").Value = 18666.111
$ws.Range("N81").Value = -20662.111

$ws.Range("H84").Value = 18049.5
$ws.Range("J84").Value = 18666.111
$ws.Range("L84").Value = 55998.333
$ws.Range("N84").Value = -65982.333

$ws.Range("H97").Value = 1645.5
$ws.Range("I97").Value = 1670.909
$ws.Range("J97").Value = 1614.4445
$ws.Range("K97").Value = 1670.909
$ws.Range("L97").Value = 1614.4445
$ws.Range("M97").Value = -1174.909
$ws.Range("N97").Value = -2606.4445

$ws.Range("H122").Value = 3897.5
$ws.Range("I122").Value = 3756.3635
$ws.Range("K122").Value = 11269.0905
$ws.Range("M122").Value = -8819.0905

$ws.Range("H132").Value = 3614.5557
$ws.Range("I132").Value = 3459.25
$ws.Range("K132").Value = 10377.75
$ws.Range("M132").Value = -7847.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4012.625
$ws.Range("J55").Value = 4129.8
$ws.Range("L55").Value = 4129.8
$ws.Range("N55").Value = -4475.8

$ws.Range("H132").Value = 4750
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

$ws.Range("H136").Value = 5834.8823
$ws.Range("I136").Value = 6242.375
$ws.Range("J136").Value = 5472.6665
$ws.Range("K136").Value = 18727.125
$ws.Range("L136").Value = 16417.9995
$ws.Range("M136").Value = -16177.125
$ws.Range("N136").Value = -21517.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 670
$ws.Range("I100").Value = 698.5
$ws.Range("K100").Value = 1397
$ws.Range("M100").Value = -856

$ws.Range("H132").Value = 52557.43
$ws.Range("I132").Value = 60567
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 181701
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -179171
$ws.Range("N132").Value = -18560
